$d = $word.ActiveDocument

# 1. Change the document's font from Garamond to Times New Roman everywhere
#    (iterate per-paragraph so the paragraph-mark run properties are
#    updated too, not just the visible run text).
foreach ($p in $d.Paragraphs) {
    $p.Range.Font.Name = "Times New Roman"
}

# 2. Remove the trailing empty paragraph before the final section break.
$count = $d.Paragraphs.Count
$last = $d.Paragraphs($count)
$prev = $d.Paragraphs($count - 1)
$rng = $d.Range($prev.Range.End - 1, $last.Range.End)
$rng.Delete()
